# Weekly crime data update (new crime data collected)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report header strings ---
$ws.Range("A8").Value = "Volume 30   Number  31"
$ws.Range("C9").Value = "Report Covering the Week  7/31/2023  Through  8/6/2023"

# --- Row 30 cells change from text placeholders ("0"/"***.*") to real numbers; set number format explicitly ---
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F30").NumberFormat = "#,##0"

# --- Update weekly crime statistics table (rows 14-30, cols C:N) ---
# Row 14
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 100
$ws.Range("F14").Value = 6
$ws.Range("G14").Value = 14
$ws.Range("H14").Value = -57.142857142857
$ws.Range("I14").Value = 77
$ws.Range("J14").Value = 89
$ws.Range("K14").Value = -13.483146067415
$ws.Range("L14").Value = -7.228915662650
$ws.Range("M14").Value = -3.75
$ws.Range("N14").Value = -73.986486486486
# Row 15
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 9
$ws.Range("E15").Value = -44.444444444444
$ws.Range("I15").Value = 230
$ws.Range("J15").Value = 246
$ws.Range("K15").Value = -6.504065040650
$ws.Range("L15").Value = 8.490566037735
$ws.Range("M15").Value = 26.373626373626
$ws.Range("N15").Value = -46.386946386946
# Row 16
$ws.Range("C16").Value = 115
$ws.Range("D16").Value = 114
$ws.Range("E16").Value = 0.877192982456
$ws.Range("F16").Value = 439
$ws.Range("H16").Value = -16.539923954372
$ws.Range("I16").Value = 2832
$ws.Range("J16").Value = 3031
$ws.Range("K16").Value = -6.565489937314
$ws.Range("L16").Value = 33.774208786018
$ws.Range("M16").Value = 9.980582524271
$ws.Range("N16").Value = -70.308240721325
# Row 17
$ws.Range("C17").Value = 156
$ws.Range("D17").Value = 133
$ws.Range("E17").Value = 17.293233082706
$ws.Range("F17").Value = 702
$ws.Range("G17").Value = 648
$ws.Range("H17").Value = 8.333333333333
$ws.Range("I17").Value = 4779
$ws.Range("J17").Value = 4405
$ws.Range("K17").Value = 8.490351872871
$ws.Range("L17").Value = 33.678321678321
$ws.Range("M17").Value = 79.391891891891
$ws.Range("N17").Value = -13.172238372093
# Row 18
$ws.Range("C18").Value = 56
$ws.Range("D18").Value = 37
$ws.Range("E18").Value = 51.351351351351
$ws.Range("F18").Value = 236
$ws.Range("G18").Value = 206
$ws.Range("H18").Value = 14.563106796116
$ws.Range("I18").Value = 1787
$ws.Range("J18").Value = 1758
$ws.Range("K18").Value = 1.649601820250
$ws.Range("L18").Value = 41.264822134387
$ws.Range("M18").Value = -5.897840968931
$ws.Range("N18").Value = -84.104251912471
# Row 19
$ws.Range("C19").Value = 183
$ws.Range("D19").Value = 154
$ws.Range("E19").Value = 18.831168831168
$ws.Range("F19").Value = 632
$ws.Range("G19").Value = 655
$ws.Range("H19").Value = -3.511450381679
$ws.Range("I19").Value = 4611
$ws.Range("J19").Value = 4739
$ws.Range("K19").Value = -2.700991770415
$ws.Range("L19").Value = 24.352750809061
$ws.Range("M19").Value = 70.273264401772
$ws.Range("N19").Value = 5.082041932543
# Row 20
$ws.Range("C20").Value = 107
$ws.Range("D20").Value = 79
$ws.Range("E20").Value = 35.443037974683
$ws.Range("F20").Value = 453
$ws.Range("G20").Value = 290
$ws.Range("H20").Value = 56.206896551724
$ws.Range("I20").Value = 3204
$ws.Range("J20").Value = 2369
$ws.Range("K20").Value = 35.246939636977
$ws.Range("L20").Value = 103.557814485388
$ws.Range("M20").Value = 160.064935064935
$ws.Range("N20").Value = -64.956797550038
# Row 21
$ws.Range("C21").Value = 624
$ws.Range("D21").Value = 527
$ws.Range("E21").Value = 18.406072106261
$ws.Range("F21").Value = 2494
$ws.Range("G21").Value = 2375
$ws.Range("H21").Value = 5.010526315789
$ws.Range("I21").Value = 17520
$ws.Range("J21").Value = 16637
$ws.Range("K21").Value = 5.307447256115
$ws.Range("L21").Value = 39.779798946864
$ws.Range("M21").Value = 54.497354497354
$ws.Range("N21").Value = -56.783423778983
# Row 22
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 18
$ws.Range("G22").Value = 21
$ws.Range("H22").Value = -14.285714285714
$ws.Range("I22").Value = 169
$ws.Range("J22").Value = 212
$ws.Range("K22").Value = -20.283018867924
$ws.Range("L22").Value = 22.463768115942
$ws.Range("M22").Value = -13.775510204081
# Row 23
$ws.Range("C23").Value = 35
$ws.Range("D23").Value = 20
$ws.Range("E23").Value = 75
$ws.Range("F23").Value = 143
$ws.Range("G23").Value = 116
$ws.Range("H23").Value = 23.275862068965
$ws.Range("I23").Value = 1061
$ws.Range("J23").Value = 959
$ws.Range("K23").Value = 10.636079249217
$ws.Range("L23").Value = 50.710227272727
$ws.Range("M23").Value = 66.040688575899
# Row 24
$ws.Range("C24").Value = 305
$ws.Range("D24").Value = 400
$ws.Range("E24").Value = -23.75
$ws.Range("F24").Value = 1431
$ws.Range("G24").Value = 1492
$ws.Range("H24").Value = -4.088471849865
$ws.Range("I24").Value = 10661
$ws.Range("J24").Value = 11004
$ws.Range("K24").Value = -3.117048346055
$ws.Range("L24").Value = 44.654002713704
$ws.Range("M24").Value = 41.693248272195
# Row 25
$ws.Range("C25").Value = 208
$ws.Range("D25").Value = 201
$ws.Range("E25").Value = 3.482587064676
$ws.Range("F25").Value = 851
$ws.Range("G25").Value = 797
$ws.Range("H25").Value = 6.775407779171
$ws.Range("I25").Value = 6322
$ws.Range("J25").Value = 6080
$ws.Range("K25").Value = 3.980263157894
$ws.Range("L25").Value = 26.896828582898
$ws.Range("M25").Value = -5.599522174107
# Row 26
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = -7.692307692307
$ws.Range("F26").Value = 47
$ws.Range("G26").Value = 51
$ws.Range("H26").Value = -7.843137254901
$ws.Range("I26").Value = 386
$ws.Range("J26").Value = 422
$ws.Range("K26").Value = -8.530805687203
$ws.Range("L26").Value = 12.209302325581
# Row 27
$ws.Range("C27").Value = 11
$ws.Range("D27").Value = 18
$ws.Range("E27").Value = -38.888888888888
$ws.Range("F27").Value = 73
$ws.Range("G27").Value = 67
$ws.Range("H27").Value = 8.955223880597
$ws.Range("I27").Value = 624
$ws.Range("J27").Value = 537
$ws.Range("K27").Value = 16.201117318435
$ws.Range("L27").Value = 19.540229885057
# Row 28
$ws.Range("C28").Value = 9
$ws.Range("D28").Value = 5
$ws.Range("E28").Value = 80
$ws.Range("F28").Value = 35
$ws.Range("G28").Value = 41
$ws.Range("H28").Value = -14.634146341463
$ws.Range("I28").Value = 244
$ws.Range("J28").Value = 319
$ws.Range("K28").Value = -23.510971786833
$ws.Range("L28").Value = -30.878186968838
$ws.Range("M28").Value = -11.913357400722
$ws.Range("N28").Value = -70.778443113772
# Row 29
$ws.Range("C29").Value = 8
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = 60
$ws.Range("F29").Value = 26
$ws.Range("G29").Value = 35
$ws.Range("H29").Value = -25.714285714285
$ws.Range("I29").Value = 198
$ws.Range("J29").Value = 273
$ws.Range("K29").Value = -27.472527472527
$ws.Range("L29").Value = -34.437086092715
$ws.Range("M29").Value = -15.021459227467
$ws.Range("N29").Value = -73.981603153745
# Row 30
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = -100
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = -75
$ws.Range("I30").Value = 13
$ws.Range("J30").Value = 31
$ws.Range("K30").Value = -58.064516129032
$ws.Range("L30").Value = -56.666666666666
